# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates (and a couple of cell deletions) across 8 sheets
# per the scraped OOXML diff for Ultros_Profits.xlsx.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1526.8
$ws.Range("I15").Value = 1526.8
$ws.Range("K15").Value = 4580.4
$ws.Range("M15").Value = -4411.4
$ws.Range("H63").Value = 150000
$ws.Range("J63").Value = 150000
$ws.Range("L63").Value = 150000
$ws.Range("N63").Value = -151248
$ws.Range("H66").Value = 150000
$ws.Range("J66").Value = 150000
$ws.Range("L66").Value = 450000
$ws.Range("N66").Value = -456240
$ws.Range("H125").Value = 6552.737
$ws.Range("I125").Value = 992.8182
$ws.Range("J125").Value = 14197.625
$ws.Range("K125").Value = 8935.363800000001
$ws.Range("L125").Value = 127778.625
$ws.Range("M125").Value = -6475.363800000001
$ws.Range("N125").Value = -132698.625
$ws.Range("H137").Value = 4129.7812
$ws.Range("I137").Value = 3476.796
$ws.Range("J137").Value = 6262.8667
$ws.Range("K137").Value = 10430.388
$ws.Range("L137").Value = 18788.6001
$ws.Range("M137").Value = -7880.387999999999
$ws.Range("N137").Value = -23888.6001
$ws.Range("H138").Value = 2332.4546
$ws.Range("J138").Value = 3689.1785
$ws.Range("L138").Value = 11067.5355
$ws.Range("N138").Value = -21347.5355

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11368492
$ws.Range("I32").Value = 12199578
$ws.Range("K32").Value = 12199578
$ws.Range("M32").Value = -12199291
$ws.Range("H61").Value = 2159.5208
$ws.Range("I61").Value = 1903.5111
$ws.Range("J61").Value = 5999.6665
$ws.Range("K61").Value = 1903.5111
$ws.Range("L61").Value = 5999.6665
$ws.Range("M61").Value = -1691.5111
$ws.Range("N61").Value = -6423.6665
$ws.Range("H136").Value = 2159.5208
$ws.Range("I136").Value = 1903.5111
$ws.Range("J136").Value = 5999.6665
$ws.Range("K136").Value = 5710.5333
$ws.Range("L136").Value = 17998.9995
$ws.Range("M136").Value = -3160.5333
$ws.Range("N136").Value = -23098.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1453.8
$ws.Range("I20").Value = 1216.25
$ws.Range("J20").Value = 1725.2858
$ws.Range("K20").Value = 1216.25
$ws.Range("L20").Value = 1725.2858
$ws.Range("M20").Value = -969.25
$ws.Range("N20").Value = -2219.2858

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 10000
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10280
$ws.Range("H22").Value = 384992.3
$ws.Range("I22").Value = 415
$ws.Range("J22").Value = 1666916.6
$ws.Range("K22").Value = 415
$ws.Range("L22").Value = 1666916.6
$ws.Range("M22").Value = -65
$ws.Range("N22").Value = -1667616.6
$ws.Range("H31").Value = 2839.606
$ws.Range("J31").Value = 2883.8572
$ws.Range("L31").Value = 2883.8572
$ws.Range("N31").Value = -3473.8572
$ws.Range("H34").Value = 2839.606
$ws.Range("J34").Value = 2883.8572
$ws.Range("L34").Value = 2883.8572
$ws.Range("N34").Value = -3287.8572
$ws.Range("H58").Value = 2865.7585
$ws.Range("I58").Value = 1288.15
$ws.Range("K58").Value = 1288.15
$ws.Range("M58").Value = -1085.15
$ws.Range("H132").Value = 1424.3334
$ws.Range("I132").Value = 1377.3334
$ws.Range("K132").Value = 4132.0002
$ws.Range("M132").Value = -1602.0002
$ws.Range("H134").Value = 3139.1887
$ws.Range("I134").Value = 2275.2683
$ws.Range("K134").Value = 6825.804900000001
$ws.Range("M134").Value = -4290.804900000001
$ws.Range("H136").Value = 2865.7585
$ws.Range("I136").Value = 1288.15
$ws.Range("K136").Value = 3864.45
$ws.Range("M136").Value = -1314.45

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 210.83333
$ws.Range("I15").Value = 62.5
$ws.Range("J15").Value = 285
$ws.Range("K15").Value = 187.5
$ws.Range("L15").Value = 855
$ws.Range("M15").Value = -47.5
$ws.Range("N15").Value = -1135
$ws.Range("H17").Value = 494.6
$ws.Range("I17").Value = 570.25
$ws.Range("J17").Value = 192
$ws.Range("K17").Value = 1710.75
$ws.Range("L17").Value = 576
$ws.Range("M17").Value = -1541.75
$ws.Range("N17").Value = -914
$ws.Range("H50").Value = 1147.5
$ws.Range("I50").Value = 200
$ws.Range("J50").Value = 1337
$ws.Range("K50").Value = 600
$ws.Range("L50").Value = 4011
$ws.Range("M50").Value = -119
$ws.Range("N50").Value = -4973
$ws.Range("H53").Value = 1147.5
$ws.Range("I53").Value = 200
$ws.Range("J53").Value = 1337
$ws.Range("K53").Value = 600
$ws.Range("L53").Value = 4011
$ws.Range("M53").Value = -119
$ws.Range("N53").Value = -4973
$ws.Range("H120").Value = 5221.4287
$ws.Range("I120").Value = 5221.4287
$ws.Range("K120").Value = 15664.2861
$ws.Range("M120").Value = -10826.2861
$ws.Range("H121").Value = 1314.25
$ws.Range("I121").Value = 307.5
$ws.Range("J121").Value = 1918.3
$ws.Range("K121").Value = 922.5
$ws.Range("L121").Value = 5754.9
$ws.Range("M121").Value = 387.5
$ws.Range("N121").Value = -8374.9

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4166898.5
$ws.Range("I2").Value = 146.16667
$ws.Range("K2").Value = 146.16667
$ws.Range("M2").Value = -33.16667000000001
$ws.Range("H11").Value = 1754785.8
$ws.Range("I11").Value = 1222091
$ws.Range("J11").Value = 3708000
$ws.Range("K11").Value = 1222091
$ws.Range("L11").Value = 3708000
$ws.Range("M11").Value = -1221952
$ws.Range("N11").Value = -3708278
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("H70").Value = 60999.57
$ws.Range("I70").Value = 116988.4
$ws.Range("J70").Value = 10100.637
$ws.Range("K70").Value = 116988.4
$ws.Range("L70").Value = 10100.637
$ws.Range("M70").Value = -116718.4
$ws.Range("N70").Value = -10640.637
$ws.Range("H73").Value = 60999.57
$ws.Range("I73").Value = 116988.4
$ws.Range("J73").Value = 10100.637
$ws.Range("K73").Value = 116988.4
$ws.Range("L73").Value = 10100.637
$ws.Range("M73").Value = -116052.4
$ws.Range("N73").Value = -11972.637
$ws.Range("H126").Value = 3690.4285
$ws.Range("I126").Value = 3690.4285
$ws.Range("K126").Value = 11071.2855
$ws.Range("M126").Value = -8601.2855
$ws.Range("H132").Value = 2862.9062
$ws.Range("I132").Value = 2623.44
$ws.Range("K132").Value = 7870.32
$ws.Range("M132").Value = -5340.32
$ws.Range("N23").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 749.5
$ws.Range("I46").Value = 499
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 499
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -311
$ws.Range("N46").Value = -1376
$ws.Range("H136").Value = 2578.3845
$ws.Range("I136").Value = 1900.8889
$ws.Range("K136").Value = 5702.6667
$ws.Range("M136").Value = -3152.6667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 495
$ws.Range("I7").Value = 495
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 495
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -382
$ws.Range("H122").Value = 353854.16
$ws.Range("I122").Value = 2517.4167
$ws.Range("J122").Value = 1407864.4
$ws.Range("K122").Value = 7552.250100000001
$ws.Range("L122").Value = 4223593.199999999
$ws.Range("M122").Value = -5102.250100000001
$ws.Range("N122").Value = -4228493.199999999
$ws.Range("N7").ClearContents()

